# Unit change: values in D/E/F/G for rows 5-8 (the AIC results block) on every
# year-sheet (2000..2100) are re-expressed in a unit one million times larger,
# i.e. every existing (non-zero) number is divided by 1,000,000.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $rng = $ws.Range("D5:G8")
    $vals = $rng.Value()
    for ($i = 1; $i -le 4; $i++) {
        for ($j = 1; $j -le 4; $j++) {
            $cur = $vals[$i, $j]
            if ($cur -ne 0) {
                $vals[$i, $j] = $cur / 1000000
            }
        }
    }
    $rng.Value = $vals
}
